# Bugfix: the "date" column (A2:A39) was stored as text labels like
# "1987Q4" .. "2024Q4" (one shared string per quarter-end). Replace them
# with real Excel date serials for each year's Q4 (Dec 31), formatted as
# a date-time, so downstream consumers get proper date values instead of
# opaque text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstYear = 1987
$lastYear = 2024
$row = 2

for ($year = $firstYear; $year -le $lastYear; $year++) {
    $quarterEnd = Get-Date -Year $year -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0
    $ws.Cells.Item($row, 1).Value = $quarterEnd.ToOADate()
    $row++
}

$ws.Range("A2:A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
